# Rotate the observation-record data held in rows 27-30 of the "Artfynd"
# sheet: each row takes on the species/observation data that previously
# belonged to the next row (row 30 wraps around to what row 27 had).
# Columns that are identical across all four rows (D, I, P, S, T, U, V, W,
# Y, AA, AD, AE, AG, AT, AW, AX, AY) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Artfynd")

# Row 27 <- old row 28 data
$ws.Range("A27").Value = 130961461
$ws.Range("B27").Value = 79245
$ws.Range("E27").Value = 6425
$ws.Range("F27").Value = "Garnlav"
$ws.Range("G27").Value = "Alectoria sarmentosa"
$ws.Range("H27").Value = "(Ach.) Ach."
$ws.Range("Q27").Value = 446088
$ws.Range("R27").Value = 6760088

# Row 28 <- old row 29 data
$ws.Range("A28").Value = 130961750
$ws.Range("Q28").Value = 446098
$ws.Range("R28").Value = 6760061
$ws.Range("AC28").Value = "Rikligt i en radie av ca 50 meter"

# Row 29 <- old row 30 data
$ws.Range("A29").Value = 130963807
$ws.Range("B29").Value = 57881
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "färska spår"
$ws.Range("Q29").Value = 445932
$ws.Range("R29").Value = 6760079
$ws.Range("Z29").Value = "14:08"
$ws.Range("AB29").Value = "14:08"
$ws.Range("AC29").ClearContents()

# Row 30 <- old row 27 data
$ws.Range("A30").Value = 130962736
$ws.Range("B30").Value = 79835
$ws.Range("E30").Value = 229821
$ws.Range("F30").Value = "Vedflamlav"
$ws.Range("G30").Value = "Ramboldia elabens"
$ws.Range("H30").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("M30").ClearContents()
$ws.Range("Q30").Value = 446008
$ws.Range("R30").Value = 6759948
$ws.Range("Z30").Value = "10:26"
$ws.Range("AB30").Value = "10:26"
